$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 header cleanup: the former "unnamed: 1_level_1" (B2) and
# "unnamed: 5_level_1" (F2) placeholder labels are replaced with "total",
# matching the already-existing "total" label used in C2.
$ws.Range("B2").Value = "total"
$ws.Range("C2").Value = "total"
$ws.Range("F2").Value = "total"
